$d = $word.ActiveDocument

# Locate the word "MVC" inside the paragraph's text.
$target = $d.Content
$found = $target.Find.Execute("MVC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $target.Start
# Force Word to split this replacement into its own run (distinct from the
# surrounding, unformatted text) by toggling a direct-formatting property
# on the selection before retyping it, then clearing that same property
# back to its original (inherited) value once the text is in place.
$target.Font.Bold = $true
$target.Text = "Hexagonal"

$newRun = $d.Range($start, $start + 9)
$newRun.Font.Bold = $false
